# Update attendance ("想去人数") and min ticket price ("最低票价") figures
# on both the "展览" and "全部类型" sheets, as per the latest data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G3").Value = 45

    $ws.Range("F6").Value = 755
    $ws.Range("F9").Value = 4505
    $ws.Range("F11").Value = 354
    $ws.Range("F12").Value = 1280
    $ws.Range("F13").Value = 531
    $ws.Range("F14").Value = 52
    $ws.Range("F15").Value = 867
    $ws.Range("F17").Value = 477
    $ws.Range("F19").Value = 230
    $ws.Range("F20").Value = 21
}
